$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price (D) / volume-change (E) values scraped for the updated symbol list.
# Values are written as text (matching the source data which stores numbers/percents
# as formatted strings), so we force Text number format before assigning, then reset
# the cell style back to Normal to avoid leaving a stray text-format style behind.
$updates = @(
    @{Cell="D2"; Value="331.30"}
    @{Cell="E2"; Value="0.08%"}
    @{Cell="D3"; Value="41.74"}
    @{Cell="E3"; Value="6.31%"}
    @{Cell="D4"; Value="5.704"}
    @{Cell="E4"; Value="0.10%"}
    @{Cell="D5"; Value="0.08351"}
    @{Cell="E5"; Value="3.83%"}
    @{Cell="D6"; Value="2.029"}
    @{Cell="E6"; Value="3.90%"}
    @{Cell="D7"; Value="8.792"}
    @{Cell="E7"; Value="2.04%"}
    @{Cell="D8"; Value="4.543"}
    @{Cell="E8"; Value="1.34%"}
    @{Cell="E9"; Value="0.59%"}
    @{Cell="D10"; Value="0.9257"}
    @{Cell="E10"; Value="0.46%"}
    @{Cell="D11"; Value="0.1293"}
    @{Cell="E11"; Value="4.48%"}
    @{Cell="D12"; Value="0.1978"}
    @{Cell="E12"; Value="1.57%"}
    @{Cell="D13"; Value="0.09510"}
    @{Cell="E13"; Value="3.43%"}
    @{Cell="D14"; Value="0.03932"}
    @{Cell="E14"; Value="12.63%"}
    @{Cell="D15"; Value="0.1059"}
    @{Cell="E15"; Value="0.93%"}
    @{Cell="D16"; Value="0.001307"}
    @{Cell="E16"; Value="0.97%"}
    @{Cell="D17"; Value="0.006108"}
    @{Cell="E17"; Value="-3.85%"}
    @{Cell="D18"; Value="3.439"}
    @{Cell="E18"; Value="2.27%"}
    @{Cell="E19"; Value="2.23%"}
    @{Cell="D20"; Value="8.233"}
    @{Cell="E20"; Value="-5.56%"}
    @{Cell="D21"; Value="0.1372"}
    @{Cell="E21"; Value="0.07%"}
    @{Cell="D22"; Value="0.2411"}
    @{Cell="E22"; Value="-7.71%"}
    @{Cell="D23"; Value="0.04436"}
    @{Cell="E23"; Value="0.23%"}
    @{Cell="E24"; Value="-0.30%"}
    @{Cell="D25"; Value="0.004372"}
    @{Cell="E25"; Value="-2.79%"}
    @{Cell="D26"; Value="0.0001201"}
    @{Cell="E26"; Value="0.06%"}
    @{Cell="D39"; Value="0.02812"}
    @{Cell="E39"; Value="9.27%"}
    @{Cell="D40"; Value="0.05548"}
    @{Cell="E40"; Value="1.95%"}
    @{Cell="D41"; Value="0.007804"}
    @{Cell="E41"; Value="3.78%"}
    @{Cell="D42"; Value="0.1436"}
    @{Cell="E42"; Value="2.41%"}
    @{Cell="D43"; Value="0.009094"}
    @{Cell="E43"; Value="-8.15%"}
    @{Cell="D44"; Value="0.002142"}
    @{Cell="E44"; Value="1.48%"}
    @{Cell="D45"; Value="0.01088"}
    @{Cell="E45"; Value="-4.90%"}
    @{Cell="D46"; Value="0.00007051"}
    @{Cell="E46"; Value="3.53%"}
    @{Cell="E47"; Value="-0.03%"}
    @{Cell="D48"; Value="0.003497"}
    @{Cell="E48"; Value="14.60%"}
    @{Cell="D49"; Value="0.002279"}
    @{Cell="E49"; Value="-0.10%"}
    @{Cell="D50"; Value="0.00002101"}
    @{Cell="E50"; Value="-0.03%"}
    @{Cell="D51"; Value="0.0002001"}
    @{Cell="E51"; Value="-0.03%"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
